$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.906.61"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "3.496.94"
$ws.Range("E3").Value = "  -2.59%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'585.03"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").Value = "'132.20"
$ws.Range("E6").Value = "  -4.20%  "
$ws.Range("D7").Value = "3.496.11"
$ws.Range("E7").Value = "  -2.70%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "'7.10"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "4.093.49"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "'27.73"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D16").Value = "'0.118"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "3.508.27"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").Value = "64.012.35"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("D19").Value = "'9.97"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").Value = "'14.52"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "'5.65"
$ws.Range("E21").Value = "  -4.00%  "
$ws.Range("D22").Value = "'390.59"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").Value = "'0.577"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").Value = "3.639.40"
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("D25").Value = "'72.99"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -4.73%  "
$ws.Range("D28").Value = "'1.58"
$ws.Range("E28").Value = "  -3.05%  "
$ws.Range("E29").Value = "  -8.30%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'2.26"
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("D32").Value = "'8.23"
$ws.Range("E32").Value = "  -4.49%  "
$ws.Range("D33").Value = "3.504.58"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'23.82"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("D38").Value = "'6.97"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -4.72%  "
$ws.Range("D40").Value = "'167.32"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Value = "'0.0808"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("D42").Value = "'27.28"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'41.96"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("D47").Value = "'4.38"
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "2.435.25"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "'0.894"
$ws.Range("E51").Value = "  -1.44%  "
